$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Rusia
$ws.Range("A10").Value = "Rusia"
$ws.Range("B10").Value = 134687
$ws.Range("C10").Value = 10633
$ws.Range("D10").Value = 16639
$ws.Range("E10").Value = 116768
$ws.Range("F10").Value = 2300
$ws.Range("G10").Value = 58
$ws.Range("H10").Value = 1280

# Row 11: Turquia
$ws.Range("A11").Value = "Turquia"
$ws.Range("B11").Value = 124375
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 58259
$ws.Range("E11").Value = 62780
$ws.Range("F11").Value = 1445
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 3336

# Row 66: Oman
$ws.Range("A66").Value = "Oman"
$ws.Range("B66").Value = 2568
$ws.Range("C66").Value = 85
$ws.Range("D66").Value = 750
$ws.Range("E66").Value = 1806
$ws.Range("F66").Value = 17
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 12

# Row 86: Eslovaquia
$ws.Range("A86").Value = "Eslovaquia"
$ws.Range("B86").Value = 1408
$ws.Range("C86").Value = 1
$ws.Range("D86").Value = 619
$ws.Range("E86").Value = 765
$ws.Range("F86").Value = 7
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 24

# Row 94: Letonia
$ws.Range("A94").Value = "Letonia"
$ws.Range("B94").Value = 879
$ws.Range("C94").Value = 8
$ws.Range("D94").Value = 348
$ws.Range("E94").Value = 515
$ws.Range("F94").Value = 4
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 16

# Row 103: Sri Lanka
$ws.Range("A103").Value = "Sri Lanka"
$ws.Range("B103").Value = 705
$ws.Range("C103").Value = 3
$ws.Range("D103").Value = 182
$ws.Range("E103").Value = 516
$ws.Range("F103").Value = 1
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 7

# Row 115: El Salvador
$ws.Range("A115").Value = "El Salvador"
$ws.Range("B115").Value = 490
$ws.Range("C115").Value = 44
$ws.Range("D115").Value = 154
$ws.Range("E115").Value = 325
$ws.Range("F115").Value = 3
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 11

# Row 116: Tanzania
$ws.Range("A116").Value = "Tanzania"
$ws.Range("B116").Value = 480
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 167
$ws.Range("E116").Value = 297
$ws.Range("F116").Value = 7
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 16

# Row 117: Malta
$ws.Range("A117").Value = "Malta"
$ws.Range("B117").Value = 468
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 379
$ws.Range("E117").Value = 85
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 4

# Row 118: Jamaica
$ws.Range("A118").Value = "Jamaica"
$ws.Range("B118").Value = 463
$ws.Range("C118").Value = 31
$ws.Range("D118").Value = 33
$ws.Range("E118").Value = 422
$ws.Range("F118").Value = 2
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 8

# Row 119: Jordania
$ws.Range("A119").Value = "Jordania"
$ws.Range("B119").Value = 460
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 367
$ws.Range("E119").Value = 84
$ws.Range("F119").Value = 5
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 9

# Row 194: Namibia
$ws.Range("A194").Value = "Namibia"
$ws.Range("B194").Value = 16
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 8
$ws.Range("E194").Value = 8
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0

# Row 195: San Vicente y las Granadinas
$ws.Range("A195").Value = "San Vicente y las Granadinas"
$ws.Range("B195").Value = 16
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 8
$ws.Range("E195").Value = 8
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0

